$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range('D4').Style

$ws.Range('D2').Value = '68.497.45'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').Value = '2.645.29'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'600.06"
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').Value = "'154.86"
$ws.Range('D6').Style = $defaultStyle
$ws.Range('E6').Value = '  +3.21%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '2.644.85'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('E10').Value = '  +6.72%  '
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = "'0.0000187"
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').Value = '  +4.00%  '
$ws.Range('D16').Value = '3.125.12'
$ws.Range('D17').Value = '68.390.82'
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').Value = '2.645.93'
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('D19').Value = "'11.40"
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').Value = '  +3.12%  '
$ws.Range('D20').Value = "'368.10"
$ws.Range('D20').Style = $defaultStyle
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('E23').Value = '  +2.43%  '
$ws.Range('E24').Value = '  +2.38%  '
$ws.Range('D25').Value = "'73.47"
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = "'10.01"
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('E28').Value = '  +6.83%  '
$ws.Range('D29').Value = '2.770.66'
$ws.Range('D30').Value = "'0.999"
$ws.Range('D30').Style = $defaultStyle
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').Value = "'575.83"
$ws.Range('D31').Style = $defaultStyle
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  +4.40%  '
$ws.Range('E33').Value = '  +4.59%  '
$ws.Range('E34').Value = '  +2.43%  '
$ws.Range('D35').Value = "'0.129"
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').Value = '  +3.28%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +3.45%  '
$ws.Range('D38').Value = "'158.80"
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Value = '  +1.87%  '
$ws.Range('E39').Value = '  +4.30%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('E42').Value = '  +3.49%  '
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('D45').Value = '0.0₆0323'
$ws.Range('E45').Value = '  +14.01%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = "'40.60"
$ws.Range('D47').Style = $defaultStyle
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').Value = "'157.19"
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('D51').Value = "'21.94"
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Value = '  +2.52%  '
